# Forgot Password Scenario Completed
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Login row's Runmode flips from YES to NO
$ws.Range("C3").Value = "NO"

# New row: Complete Course scenario
$ws.Range("A5").Value = "Complete Course"
$ws.Range("B5").Value = "Complete course description"
$ws.Range("C5").Value = "NO"

# Move the active selection like the author left it
[void]$ws.Range("B14").Select()
